# Auto-generated edit script: add new date column (2025-12-08) to all data sheets
$wb = $excel.ActiveWorkbook

$sheetConfigs = @(
    @{ SheetIndex = 10; NewCol = 10; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(29,46,33,33,54,58,48,40,41,69,70,78,49,57,60,57,59,66) },
    @{ SheetIndex = 11; NewCol = 50; IsTextHeader = $True; HeaderValue = 20251208; DataValues = @(94,103,97,101,104,105,103,106,101,109,111,113,109,98,107,108,102,100) },
    @{ SheetIndex = 12; NewCol = 31; IsTextHeader = $True; HeaderValue = 20251208; DataValues = @(-36.33,-1.88,-28.52,-2.87,-11.06,-4.77,-0.04,-9.28,-3.89,3.31,1.5,10.44,-8.85,-26.63,-6.98,-2.85,-4.83,-25.68) },
    @{ SheetIndex = 13; NewCol = 10; IsTextHeader = $True; HeaderValue = 20251208; DataValues = @(142,30,55,10,32,31,25,56,20,61,117,127,66,62,56,28,11,57) },
    @{ SheetIndex = 2; NewCol = 69; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(15690,57500,18440,16230,16905,19605,18425,1550,17210,5395,5830,5805,1689,19660,11000,6490,18225,16200) },
    @{ SheetIndex = 3; NewCol = 69; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(15700,58310,18490,16275,17260,19830,18785,1658,17360,5740,6250,6315,1818,19665,11600,6890,18430,16200) },
    @{ SheetIndex = 4; NewCol = 69; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(14905,56625,17880,15950,16692,19335,18190,1527,17000,5385,5805,5800,1660,18730,10935,6485,17975,15840) },
    @{ SheetIndex = 5; NewCol = 69; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(14985,58310,17975,16275,17260,19830,18775,1652,17315,5725,6240,6300,1818,18830,11570,6890,18420,15905) },
    @{ SheetIndex = 6; NewCol = 69; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(294880,1398294,166500,169911,1261477,2361249,3598130,8073534,600293,1503049,6275949,16200781,66152109,1995515,672002,135740,38888,54410) },
    @{ SheetIndex = 7; NewCol = 50; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(25,80,44,70,96,100,81,59,79,100,100,100,70,62,95,85,84,72) },
    @{ SheetIndex = 8; NewCol = 10; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(60,87,65,83,93,94,91,66,85,99,98,100,75,74,90,87,95,80) },
    @{ SheetIndex = 9; NewCol = 50; IsTextHeader = $False; HeaderValue = 20251208; DataValues = @(-90,44,-52,24,59,69,46,37,22,103,110,120,53,-16,86,76,37,4) }
)


foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.SheetIndex)
    $oldCol = $cfg.NewCol - 1

    # 1) Copy the column width from the previous last column to the new column
    $ws.Cells.Item(1, $cfg.NewCol).ColumnWidth = $ws.Cells.Item(1, $oldCol).ColumnWidth

    # 2) Write the header cell (row 1)
    if ($cfg.IsTextHeader) {
        # Force text storage so it matches the sibling header cells (t="inlineStr")
        $ws.Cells.Item(1, $cfg.NewCol).NumberFormat = "@"
        $ws.Cells.Item(1, $cfg.NewCol).Value = [string]$cfg.HeaderValue
    } else {
        $ws.Cells.Item(1, $cfg.NewCol).Value = $cfg.HeaderValue
    }

    # 3) Copy formatting (bold + gray fill) from the neighboring header cell onto the new one
    $ws.Cells.Item(1, $oldCol).Copy()
    $ws.Cells.Item(1, $cfg.NewCol).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # 4) Write the data rows (2 through 19)
    for ($i = 0; $i -lt $cfg.DataValues.Length; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, $cfg.NewCol).Value = $cfg.DataValues[$i]
    }
}

Write-Host "Update complete"
